$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the swapped "smallest_area_(m^3)" / "area_standard_derivation_(+-)"
#    headers (column N should read "area_standard_derivation_(+-)", column O
#    should read "smallest_area_(m^3)") and move the matching data values.
# ---------------------------------------------------------------------------
$ws.Range("N1").Value = "area_standard_derivation_(+-)"
$ws.Range("O1").Value = "smallest_area_(m^3)"

# Row 2 only had a single (avg=min=max) figure stored under N; it belongs in O.
$v2 = $ws.Range("N2").Value2()
$ws.Range("N2").ClearContents()
$ws.Range("O2").Value = $v2

# Row 3 - same situation as row 2.
$v3 = $ws.Range("N3").Value2()
$ws.Range("N3").ClearContents()
$ws.Range("O3").Value = $v3

# Rows 4-6 had both figures, just swapped between the two columns.
$n4 = $ws.Range("N4").Value2()
$o4 = $ws.Range("O4").Value2()
$ws.Range("N4").Value = $o4
$ws.Range("O4").Value = $n4

$n5 = $ws.Range("N5").Value2()
$o5 = $ws.Range("O5").Value2()
$ws.Range("N5").Value = $o5
$ws.Range("O5").Value = $n5

$n6 = $ws.Range("N6").Value2()
$o6 = $ws.Range("O6").Value2()
$ws.Range("N6").Value = $o6
$ws.Range("O6").Value = $n6

# ---------------------------------------------------------------------------
# 2. Body formatting: left/top aligned, wrapped text for every data cell
#    (build the format on A2 then fan it out over the rest of the table so
#    only one new style is created instead of one per property assignment).
# ---------------------------------------------------------------------------
$fmtSrc = $ws.Range("A2")
$fmtSrc.HorizontalAlignment = -4131   # xlLeft
$fmtSrc.VerticalAlignment = -4160     # xlTop
$fmtSrc.WrapText = $true
$fmtSrc.Copy()
$ws.Range("A2:R6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Column widths matching the real workbook's manual column sizing.
# ---------------------------------------------------------------------------
$ws.Range("A1:F1").ColumnWidth = 33.83333333333333
$ws.Range("G1").ColumnWidth = 13.333333333333332
$ws.Range("H1").ColumnWidth = 31.666666666666668
$ws.Range("I1").ColumnWidth = 17.666666666666668
$ws.Range("J1").ColumnWidth = 30.0
$ws.Range("K1").ColumnWidth = 64.33333333333334
$ws.Range("L1").ColumnWidth = 21.333333333333336
$ws.Range("M1").ColumnWidth = 18.166666666666668
$ws.Range("N1").ColumnWidth = 29.333333333333336
$ws.Range("O1").ColumnWidth = 20.833333333333336
$ws.Range("P1").ColumnWidth = 19.5
$ws.Range("Q1").ColumnWidth = 11.0
$ws.Range("R1").ColumnWidth = 14.0

# ---------------------------------------------------------------------------
# 4. Turn the header row into an AutoFilter (adds the hidden _FilterDatabase
#    defined name automatically scoped to this sheet).
# ---------------------------------------------------------------------------
$ws.Range("A1:R6").AutoFilter()
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$R`$6")
$fdb.Visible = $false

Write-Host "edit complete"
